$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new day-count entry for row 13 (B13 = 2)
$ws.Range("B13").Value = 2

# Swap the "estimate" labels between H4 and H5 (with their styles)
$ws.Range("H5").Copy()
$ws.Range("H4").Value = "1h"
$ws.Range("H5").Value = "1.5-2h"

# Update the active selection to I6
$ws.Range("I6").Select()
